$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> (DAMSLTag, DialogAct) after re-running SGNN annotation
$updates = @(
    @(33, 'sv', 'Statement-opinion'),
    @(51, 'sv', 'Statement-opinion'),
    @(59, 'sd', 'Statement-non-opinion'),
    @(61, 'sd', 'Statement-non-opinion'),
    @(62, 'sd', 'Statement-non-opinion'),
    @(63, 'sd', 'Statement-non-opinion'),
    @(64, 'sd', 'Statement-non-opinion'),
    @(65, 'sd', 'Statement-non-opinion'),
    @(115, 'sd', 'Statement-non-opinion'),
    @(117, '%', 'Uninterpretable'),
    @(132, 'sv', 'Statement-opinion'),
    @(133, 'aa', 'Agree/Accept'),
    @(179, 'aa', 'Agree/Accept'),
    @(183, '%', 'Uninterpretable'),
    @(187, 'sv', 'Statement-opinion'),
    @(219, '%', 'Uninterpretable'),
    @(237, 'sd', 'Statement-non-opinion'),
    @(254, 'sv', 'Statement-opinion'),
    @(259, 'aa', 'Agree/Accept'),
    @(262, 'sd', 'Statement-non-opinion'),
    @(268, 'sd', 'Statement-non-opinion'),
    @(273, 'sd', 'Statement-non-opinion'),
    @(274, 'sd', 'Statement-non-opinion'),
    @(280, 'sd', 'Statement-non-opinion'),
    @(281, 'sd', 'Statement-non-opinion'),
    @(283, 'sd', 'Statement-non-opinion'),
    @(289, 'sd', 'Statement-non-opinion'),
    @(293, 'sd', 'Statement-non-opinion'),
    @(297, 'sd', 'Statement-non-opinion'),
    @(308, 'sd', 'Statement-non-opinion'),
    @(309, 'sd', 'Statement-non-opinion'),
    @(314, 'aa', 'Agree/Accept'),
    @(316, 'sd', 'Statement-non-opinion'),
    @(329, 'sv', 'Statement-opinion'),
    @(332, 'sd', 'Statement-non-opinion'),
    @(359, 'sd', 'Statement-non-opinion'),
    @(364, 'sd', 'Statement-non-opinion'),
    @(371, 'sd', 'Statement-non-opinion'),
    @(383, 'sv', 'Statement-opinion'),
    @(412, 'sd', 'Statement-non-opinion'),
    @(413, 'sv', 'Statement-opinion'),
    @(419, 'sv', 'Statement-opinion'),
    @(420, 'sv', 'Statement-opinion'),
    @(448, 'sd', 'Statement-non-opinion'),
    @(462, 'aa', 'Agree/Accept'),
    @(464, 'aa', 'Agree/Accept'),
    @(465, 'sd', 'Statement-non-opinion'),
    @(466, 'sd', 'Statement-non-opinion'),
    @(468, 'sd', 'Statement-non-opinion'),
    @(470, 'sd', 'Statement-non-opinion'),
    @(471, 'sd', 'Statement-non-opinion'),
    @(472, 'sd', 'Statement-non-opinion'),
    @(486, '%', 'Uninterpretable'),
    @(487, 'sd', 'Statement-non-opinion'),
    @(510, 'sv', 'Statement-opinion'),
    @(511, 'sd', 'Statement-non-opinion'),
    @(517, 'sv', 'Statement-opinion'),
    @(522, 'sd', 'Statement-non-opinion'),
    @(525, 'sv', 'Statement-opinion'),
    @(542, 'sv', 'Statement-opinion'),
    @(549, 'sv', 'Statement-opinion'),
    @(551, '%', 'Uninterpretable'),
    @(556, 'sd', 'Statement-non-opinion'),
    @(559, 'sd', 'Statement-non-opinion'),
    @(564, '%', 'Uninterpretable'),
    @(582, '%', 'Uninterpretable'),
    @(587, 'sd', 'Statement-non-opinion'),
    @(595, '%', 'Uninterpretable'),
    @(613, 'sd', 'Statement-non-opinion'),
    @(623, '%', 'Uninterpretable'),
    @(630, '%', 'Uninterpretable'),
    @(632, '%', 'Uninterpretable'),
    @(658, 'sd', 'Statement-non-opinion'),
    @(685, 'aa', 'Agree/Accept'),
    @(687, 'sv', 'Statement-opinion'),
    @(699, 'sd', 'Statement-non-opinion'),
    @(711, 'sd', 'Statement-non-opinion'),
    @(717, 'sd', 'Statement-non-opinion'),
    @(718, 'sd', 'Statement-non-opinion'),
    @(722, 'sd', 'Statement-non-opinion'),
    @(731, 'sd', 'Statement-non-opinion'),
    @(734, 'sv', 'Statement-opinion'),
    @(736, 'aa', 'Agree/Accept'),
    @(737, 'aa', 'Agree/Accept'),
    @(778, 'sd', 'Statement-non-opinion'),
    @(783, 'sd', 'Statement-non-opinion'),
    @(793, 'sd', 'Statement-non-opinion'),
    @(817, 'sd', 'Statement-non-opinion'),
    @(830, 'sd', 'Statement-non-opinion'),
    @(846, 'sv', 'Statement-opinion'),
    @(847, 'sd', 'Statement-non-opinion'),
    @(858, 'sd', 'Statement-non-opinion'),
    @(859, 'sv', 'Statement-opinion')
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 9).Value = $u[1]
    $ws.Cells.Item($row, 10).Value = $u[2]
}
